# Updated cryptos list - applies the Price (D) and Volume(1h) (E) changes
# from the commit, including the EnergySwap/RenderToken row swap at 48/49.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Price column values look numeric ("242.84", "1.0000", ...); force
    # text storage so Excel does not coerce them into Number/Date cells,
    # then drop back to the Normal style so no stray numFmt lingers.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.820.28"
$ws.Range("E2").Value = "  -0.45%  "

Set-TextValue $ws.Range("D3") "1.895.10"
$ws.Range("E3").Value = "  +0.05%  "

Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.15%  "

Set-TextValue $ws.Range("D5") "0.7979"
$ws.Range("E5").Value = "  -2.86%  "

Set-TextValue $ws.Range("D6") "242.84"
$ws.Range("E6").Value = "  +0.44%  "

Set-TextValue $ws.Range("D7") "1.000"
$ws.Range("E7").Value = "  -0.24%  "

Set-TextValue $ws.Range("D8") "0.3169"
$ws.Range("E8").Value = "  -1.82%  "

Set-TextValue $ws.Range("D9") "25.48"
$ws.Range("E9").Value = "  -3.85%  "

Set-TextValue $ws.Range("D10") "0.07041"
$ws.Range("E10").Value = "  +0.22%  "

Set-TextValue $ws.Range("D11") "0.08060"
$ws.Range("E11").Value = "  +0.28%  "

Set-TextValue $ws.Range("D12") "0.7712"
$ws.Range("E12").Value = "  +3.14%  "

Set-TextValue $ws.Range("D13") "1.889.97"
$ws.Range("E13").Value = "  -0.15%  "

Set-TextValue $ws.Range("D14") "5.334"
$ws.Range("E14").Value = "  +2.59%  "

Set-TextValue $ws.Range("D15") "92.41"
$ws.Range("E15").Value = "  +0.09%  "

Set-TextValue $ws.Range("D16") "29.825.53"
$ws.Range("E16").Value = "  -0.32%  "

Set-TextValue $ws.Range("D17") "6.017"
$ws.Range("E17").Value = "  +2.06%  "

Set-TextValue $ws.Range("D18") "13.88"
$ws.Range("E18").Value = "  -1.09%  "

Set-TextValue $ws.Range("D19") "244.44"
$ws.Range("E19").Value = "  -0.17%  "

Set-TextValue $ws.Range("D20") "0.000007713"
$ws.Range("E20").Value = "  -0.47%  "

Set-TextValue $ws.Range("D21") "8.260"
$ws.Range("E21").Value = "  +19.60%  "

Set-TextValue $ws.Range("D22") "1.0000"
$ws.Range("E22").Value = "  -0.22%  "

Set-TextValue $ws.Range("D23") "2.136.81"
$ws.Range("E23").Value = "  -0.12%  "

Set-TextValue $ws.Range("D24") "1.001"
$ws.Range("E24").Value = "  -0.12%  "

Set-TextValue $ws.Range("D25") "0.1652"
$ws.Range("E25").Value = "  +4.08%  "

Set-TextValue $ws.Range("D26") "9.340"
$ws.Range("E26").Value = "  +1.61%  "

$ws.Range("E27").Value = "  -0.29%  "

Set-TextValue $ws.Range("D28") "18.69"
$ws.Range("E28").Value = "  -0.76%  "

Set-TextValue $ws.Range("D29") "2.058"
$ws.Range("E29").Value = "  -0.56%  "

Set-TextValue $ws.Range("D30") "1.409"
$ws.Range("E30").Value = "  +2.73%  "

$ws.Range("E31").Value = "  +1.42%  "

Set-TextValue $ws.Range("D32") "4.438"
$ws.Range("E32").Value = "  +4.05%  "

Set-TextValue $ws.Range("D33") "0.05712"
$ws.Range("E33").Value = "  +1.80%  "

Set-TextValue $ws.Range("D34") "4.047"
$ws.Range("E34").Value = "  -0.71%  "

Set-TextValue $ws.Range("D35") "1.261"
$ws.Range("E35").Value = "  -0.91%  "

Set-TextValue $ws.Range("D36") "0.7388"
$ws.Range("E36").Value = "  +1.17%  "

Set-TextValue $ws.Range("D37") "0.9983"
$ws.Range("E37").Value = "  -0.30%  "

Set-TextValue $ws.Range("D38") "2.619"
$ws.Range("E38").Value = "  -3.90%  "

Set-TextValue $ws.Range("D39") "0.01910"
$ws.Range("E39").Value = "  -0.20%  "

Set-TextValue $ws.Range("D40") "2.787"

Set-TextValue $ws.Range("D41") "0.4408"
$ws.Range("E41").Value = "  -0.18%  "

Set-TextValue $ws.Range("D42") "72.56"
$ws.Range("E42").Value = "  +0.86%  "

Set-TextValue $ws.Range("D43") "5.815"
$ws.Range("E43").Value = "  -2.37%  "

Set-TextValue $ws.Range("D44") "0.8416"
$ws.Range("E44").Value = "  -0.25%  "

Set-TextValue $ws.Range("D45") "0.9999"
$ws.Range("E45").Value = "  -0.32%  "

Set-TextValue $ws.Range("D46") "1.033.59"
$ws.Range("E46").Value = "  +4.33%  "

Set-TextValue $ws.Range("D47") "102.80"
$ws.Range("E47").Value = "  +2.11%  "

Set-TextValue $ws.Range("D50") "7.432"
$ws.Range("E50").Value = "  -1.97%  "

Set-TextValue $ws.Range("D51") "2.035.25"
$ws.Range("E51").Value = "  -0.39%  "

# Rows 48-49: RenderToken and EnergySwap swapped places with updated data
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "10.04"
$ws.Range("E48").Value = "  +3.73%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D49") "1.873"
$ws.Range("E49").Value = "  -0.27%  "
